$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Clean whitespace / stray punctuation in team & class names (shared text) ---
# Collapses runs of 2+ spaces to a single space and strips stray apostrophes,
# applied as exact substring replacements so only the intended shared string is touched.
$ws.Cells.Replace("ATLANTIS  -  FLA", "ATLANTIS - FLA") | Out-Null
$ws.Cells.Replace("LINCOLN -  CALIFONIA", "LINCOLN - CALIFONIA") | Out-Null
$ws.Cells.Replace("BLUEFIELD  -  VA", "BLUEFIELD - VA") | Out-Null
$ws.Cells.Replace("MADONNA  - MICHIGAN", "MADONNA - MICHIGAN") | Out-Null
$ws.Cells.Replace("MID-AMERICA  NAZARENE", "MID-AMERICA NAZARENE") | Out-Null
$ws.Cells.Replace("MT.  MARTY    S.D. ", "MT. MARTY S.D. ") | Out-Null
$ws.Cells.Replace("OTTAWA  (AZ)", "OTTAWA (AZ)") | Out-Null
$ws.Cells.Replace("THOMAS  -  GA", "THOMAS - GA") | Out-Null
$ws.Cells.Replace("ALABAMA  A&M", "ALABAMA A&M") | Out-Null
$ws.Cells.Replace("DELAWARE  STATE", "DELAWARE STATE") | Out-Null
$ws.Cells.Replace("FLORIDA  A&M", "FLORIDA A&M") | Out-Null
$ws.Cells.Replace("NORTH CAROLINA  A&T", "NORTH CAROLINA A&T") | Out-Null
$ws.Cells.Replace("SOUTH EASTERN  MISSOURI", "SOUTH EASTERN MISSOURI") | Out-Null
$ws.Cells.Replace("ALLEN  (SC)", "ALLEN (SC)") | Out-Null
$ws.Cells.Replace("NEBRASKA -  KEARNEY", "NEBRASKA - KEARNEY") | Out-Null
$ws.Cells.Replace("NO'EASTERN   OKLA ST.", "NOEASTERN OKLA ST.") | Out-Null
$ws.Cells.Replace("SOUTH DAKOTA  TECH", "SOUTH DAKOTA TECH") | Out-Null
$ws.Cells.Replace("WHEELING  UNIVERSITY", "WHEELING UNIVERSITY") | Out-Null
$ws.Cells.Replace("BETHANY  -  WV", "BETHANY - WV") | Out-Null
$ws.Cells.Replace("EASTERN  PA", "EASTERN PA") | Out-Null
$ws.Cells.Replace("KING'S COLLEGE", "KINGS COLLEGE") | Out-Null
$ws.Cells.Replace("NEBRASKA  WESLEYAN", "NEBRASKA WESLEYAN") | Out-Null
$ws.Cells.Replace("U  O F  NEW ENGLAND", "U O F NEW ENGLAND") | Out-Null
$ws.Cells.Replace("HAWAI'I", "HAWAII") | Out-Null
$ws.Cells.Replace("NORTH  CAROLINA", "NORTH CAROLINA") | Out-Null
$ws.Cells.Replace("NORTHERN  ILLINOIS", "NORTHERN ILLINOIS") | Out-Null
$ws.Cells.Replace("NORTH  TEXAS", "NORTH TEXAS") | Out-Null
$ws.Cells.Replace("SOUTHERN  MISSISSIPPI", "SOUTHERN MISSISSIPPI") | Out-Null
$ws.Cells.Replace("TEXAS  STATE-SAN MARCOS", "TEXAS STATE-SAN MARCOS") | Out-Null
$ws.Cells.Replace("WEST  VIRGINIA", "WEST VIRGINIA") | Out-Null
$ws.Cells.Replace("DIVISION 1  FBS", "DIVISION 1 FBS") | Out-Null

# --- Fix mislabeled team abbreviations ---
# NWST must be remapped to WIU only *after* the original WIU has already
# become EIU, otherwise both would collapse onto the same abbreviation.
$ws.Cells.Replace("WIU", "EIU") | Out-Null
$ws.Cells.Replace("HCU", "ACU") | Out-Null
$ws.Cells.Replace("YSU", "JKST") | Out-Null
$ws.Cells.Replace("NWST", "WIU") | Out-Null

# --- Updated confidence (column F) values ---
$ws.Range("F673").Value = 83
$ws.Range("F711").Value = 86
$ws.Range("F712").Value = 88
$ws.Range("F713").Value = 82
$ws.Range("F735").Value = 50
$ws.Range("F746").Value = 47
$ws.Range("F763").Value = 56
$ws.Range("F764").Value = 85
